$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:H1 (reuse the bold/centered/bordered header style from A1:E1) ---
$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "source_file"
$ws.Range("H1").Value = "text"

$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- New data row 2 ---
$ws.Range("A2").Value = "Sunsi Wu"

# politeness_score "3" is stored as text (not a number) in the source data.
# Stage it in a scratch cell with a leading apostrophe (forces text entry),
# then paste-values-only into B2 so the text type carries over without
# leaving the quote-prefix formatting behind on B2 itself.
$ws.Range("ZZ1").Value = "'3"
$ws.Range("ZZ1").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("C2").Value = "how"
$ws.Range("D2").Value = "QSN"
$ws.Range("E2").Value = "MET"
$ws.Range("F2").Value = "1269f1fb-9c21-42a9-ae5e-c80f92622adc"
$ws.Range("G2").Value = "Bk6qQGWRb_annotated.xlsx"
$ws.Range("H2").Value = "Then how bootstrap dqn extend the idea to deep learning, followed by the noisy net, bbq, shallow UBE and LS-DQN."
